$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 'Compass Error Compass Error Compass data error Please contact DJI Support .'
$ws.Range("C2").Value = 'Compass Error'
$ws.Range("D2").Value = '0-1'
$ws.Range("E2").Value = 'Missing'
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 'Compass Error Compass Error Compass data error Please contact DJI Support .'
$ws.Range("C3").Value = 'Compass Error'
$ws.Range("D3").Value = '2-3'
$ws.Range("E3").Value = 'Missing'
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 'Compass Error Compass Error Compass data error Please contact DJI Support .'
$ws.Range("C4").Value = 'Compass data error'
$ws.Range("D4").Value = '4-6'
$ws.Range("E4").Value = 'Missing'
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 'Compass Error Compass Error Compass data error Please contact DJI Support .'
$ws.Range("C5").Value = 'Please contact DJI Support'
$ws.Range("D5").Value = '7-10'
$ws.Range("E5").Value = 'Missing'
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn .'
$ws.Range("C6").Value = 'GPS signal weak'
$ws.Range("D6").Value = '0-2'
$ws.Range("E6").Value = 'Missing'
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 'GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn .'
$ws.Range("C7").Value = 'Fly with caution'
$ws.Range("D7").Value = '3-5'
$ws.Range("E7").Value = 'Missing'
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 'GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn .'
$ws.Range("C8").Value = 'Aircraft in Altitude Zone'
$ws.Range("D8").Value = '6-9'
$ws.Range("E8").Value = 'Missing'
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 'GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn .'
$ws.Range("C9").Value = 'Max altitude set to nnn'
$ws.Range("D9").Value = '10-14'
$ws.Range("E9").Value = 'Missing'
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 'GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn .'
$ws.Range("C10").Value = 'GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn'
$ws.Range("D10").Value = '0-14'
$ws.Range("E10").Value = "'False"
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 'High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP .'
$ws.Range("C11").Value = 'High wind velocity'
$ws.Range("D11").Value = '0-2'
$ws.Range("E11").Value = 'Missing'
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = 'High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP .'
$ws.Range("C12").Value = 'Ensure the aircraft remains within your line of sight and fly with caution'
$ws.Range("D12").Value = '3-15'
$ws.Range("E12").Value = 'Missing'
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 'High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP .'
$ws.Range("C13").Value = 'High wind velocity'
$ws.Range("D13").Value = '16-18'
$ws.Range("E13").Value = 'Missing'
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = 'High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP .'
$ws.Range("C14").Value = 'Fly with caution and land in a safe place ASAP'
$ws.Range("D14").Value = '19-28'
$ws.Range("E14").Value = 'Missing'
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 'High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP .'
$ws.Range("C15").Value = 'High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP'
$ws.Range("D15").Value = '0-28'
$ws.Range("E15").Value = "'False"
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 'Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support .'
$ws.Range("C16").Value = 'Motor speed error'
$ws.Range("D16").Value = '0-2'
$ws.Range("E16").Value = 'Missing'
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = 'Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support .'
$ws.Range("C17").Value = 'Land or return to home promptly'
$ws.Range("D17").Value = '3-8'
$ws.Range("E17").Value = 'Missing'
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = 'Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support .'
$ws.Range("C18").Value = 'After powering off the aircraft, replace the propeller on the beeping ESC'
$ws.Range("D18").Value = '9-20'
$ws.Range("E18").Value = 'Missing'
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = 'Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support .'
$ws.Range("C19").Value = 'If the issue persists, contact DJI Support'
$ws.Range("D19").Value = '21-27'
$ws.Range("E19").Value = 'Missing'
